$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.638.69"
$ws.Range("E2").Value = "  -2.27%  "
$ws.Range("D3").Value = "2.418.97"
$ws.Range("E3").Value = "  +6.93%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "293.98"
$ws.Range("E5").Value = "  -1.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.56"
$ws.Range("E6").Value = "  -5.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.567"
$ws.Range("E7").Value = "  +1.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("E9").Value = "  -1.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.78"
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("E11").Value = "  -3.69%  "
$ws.Range("E12").Value = "  -0.34%  "
$ws.Range("E13").Value = "  +2.18%  "
$ws.Range("D14").Value = "2.787.32"
$ws.Range("E14").Value = "  +6.78%  "
$ws.Range("D15").Value = "2.412.00"
$ws.Range("E15").Value = "  +6.62%  "
$ws.Range("E16").Value = "  +6.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.13"
$ws.Range("E17").Value = "  +3.80%  "
$ws.Range("D18").Value = "45.510.94"
$ws.Range("E18").Value = "  -2.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.30"
$ws.Range("E19").Value = "  -3.57%  "
$ws.Range("E20").Value = "  -1.37%  "
$ws.Range("E21").Value = "  +6.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.12"
$ws.Range("E22").Value = "  +2.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "242.06"
$ws.Range("E23").Value = "  -1.86%  "
$ws.Range("E24").Value = "  -0.47%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  +4.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "38.70"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.22"
$ws.Range("E28").Value = "  -1.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.69"
$ws.Range("E29").Value = "  +0.91%  "
$ws.Range("B30").Value = "LidoDAOToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.81"
$ws.Range("E30").Value = "  +15.70%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "21.19"
$ws.Range("E31").Value = "  +3.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.77"
$ws.Range("E32").Value = "  -2.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "147.64"
$ws.Range("E33").Value = "  +1.15%  "
$ws.Range("E34").Value = "  +1.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0765"
$ws.Range("E35").Value = "  -0.25%  "
$ws.Range("E36").Value = "  +17.84%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.115"
$ws.Range("E38").Value = "  +0.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.56"
$ws.Range("E39").Value = "  -6.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.78"
$ws.Range("E40").Value = "  -1.53%  "
$ws.Range("E41").Value = "  +0.44%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.988.32"
$ws.Range("E42").Value = "  +11.33%  "
$ws.Range("B43").Value = "NEARProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.20"
$ws.Range("E43").Value = "  +3.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "89.28"
$ws.Range("E45").Value = "  -1.55%  "
$ws.Range("E46").Value = "  -2.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.11"
$ws.Range("E47").Value = "  +26.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.53"
$ws.Range("E48").Value = "  +9.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "100.64"
$ws.Range("E49").Value = "  +6.84%  "
$ws.Range("D50").Value = "2.656.15"
$ws.Range("E50").Value = "  +6.77%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.183"
$ws.Range("E51").Value = "  -1.08%  "
